# Update "想去人数" (interested-count) figures, and a couple of
# "最低票价" (min ticket price) corrections, across all four sheets
# of the workbook, as published by the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

function Set-Cell {
    param([string]$SheetName, [string]$CellRef, [double]$NewValue)
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($CellRef).Value = $NewValue
}

# ---- 展览 (Exhibitions) ----
Set-Cell "展览" "F2"  68
Set-Cell "展览" "F3"  393
Set-Cell "展览" "F4"  167
Set-Cell "展览" "F5"  1337
Set-Cell "展览" "F6"  238
Set-Cell "展览" "F7"  2538
Set-Cell "展览" "F8"  941
Set-Cell "展览" "F9"  18861
Set-Cell "展览" "F10" 57
Set-Cell "展览" "F11" 1981
Set-Cell "展览" "F12" 681
Set-Cell "展览" "G12" 54
Set-Cell "展览" "F13" 604
Set-Cell "展览" "F14" 346
Set-Cell "展览" "F15" 616
Set-Cell "展览" "G15" 54
Set-Cell "展览" "F16" 201
Set-Cell "展览" "F19" 328
Set-Cell "展览" "F20" 46
Set-Cell "展览" "F21" 210
Set-Cell "展览" "F23" 119

# ---- 演出 (Performances) ----
Set-Cell "演出" "F5"  175
Set-Cell "演出" "F10" 235
Set-Cell "演出" "F16" 75
Set-Cell "演出" "F18" 19

# ---- 本地生活 (Local Life) ----
Set-Cell "本地生活" "F2" 5911
Set-Cell "本地生活" "F3" 584

# ---- 全部类型 (All Types - aggregated view) ----
Set-Cell "全部类型" "F2"  68
Set-Cell "全部类型" "F3"  584
Set-Cell "全部类型" "F5"  393
Set-Cell "全部类型" "F7"  167
Set-Cell "全部类型" "F9"  1337
Set-Cell "全部类型" "F11" 238
Set-Cell "全部类型" "F12" 175
Set-Cell "全部类型" "F14" 2538
Set-Cell "全部类型" "F15" 941
Set-Cell "全部类型" "F16" 18861
Set-Cell "全部类型" "F19" 57
Set-Cell "全部类型" "F21" 235
Set-Cell "全部类型" "F22" 1981
Set-Cell "全部类型" "F23" 681
Set-Cell "全部类型" "G23" 54
Set-Cell "全部类型" "F25" 346
Set-Cell "全部类型" "F26" 616
Set-Cell "全部类型" "G26" 54
Set-Cell "全部类型" "F27" 201
Set-Cell "全部类型" "F32" 328
Set-Cell "全部类型" "F33" 46
Set-Cell "全部类型" "F35" 210
Set-Cell "全部类型" "F36" 75
Set-Cell "全部类型" "F38" 119
Set-Cell "全部类型" "F39" 19
